$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.737.05"
$ws.Range("E2").Value = "  +0.17%  "
$ws.Range("D3").Value = "'2.286.52"
$ws.Range("E3").Value = "  -0.16%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.38%  "
$ws.Range("D5").Value = "'110.70"
$ws.Range("E5").Value = "  +15.39%  "
$ws.Range("D6").Value = "'268.14"
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("E7").Value = "  +0.54%  "
$ws.Range("E8").Value = "  +0.31%  "
$ws.Range("D9").Value = "'0.618"
$ws.Range("E9").Value = "  +1.52%  "
$ws.Range("D10").Value = "'47.72"
$ws.Range("E10").Value = "  +5.11%  "
$ws.Range("E11").Value = "  +1.54%  "
$ws.Range("D12").Value = "'9.10"
$ws.Range("E12").Value = "  +15.06%  "
$ws.Range("D13").Value = "'0.106"
$ws.Range("E13").Value = "  +0.04%  "
$ws.Range("D14").Value = "'15.82"
$ws.Range("E14").Value = "  +3.11%  "
$ws.Range("D15").Value = "'2.630.24"
$ws.Range("E15").Value = "  -0.06%  "
$ws.Range("E16").Value = "  -0.11%  "
$ws.Range("D17").Value = "'2.278.81"
$ws.Range("E17").Value = "  -0.29%  "
$ws.Range("D18").Value = "'43.619.49"
$ws.Range("E18").Value = "  +0.11%  "
$ws.Range("E19").Value = "  +0.51%  "
$ws.Range("D20").Value = "'6.74"
$ws.Range("E20").Value = "  +8.74%  "
$ws.Range("D21").Value = "'72.11"
$ws.Range("E21").Value = "  +0.06%  "
$ws.Range("E22").Value = "  -4.59%  "
$ws.Range("D23").Value = "'9.81"
$ws.Range("E23").Value = "  +7.44%  "
$ws.Range("D24").Value = "'232.37"
$ws.Range("E24").Value = "  -0.08%  "
$ws.Range("D25").Value = "'2.78"
$ws.Range("E25").Value = "  +6.98%  "
$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").Value = "'11.69"
$ws.Range("E27").Value = "  +4.22%  "
$ws.Range("D28").Value = "'41.69"
$ws.Range("E28").Value = "  +4.02%  "
$ws.Range("E29").Value = "  -2.22%  "
$ws.Range("E30").Value = "  -0.04%  "
$ws.Range("D31").Value = "'175.84"
$ws.Range("E31").Value = "  +0.55%  "
$ws.Range("D32").Value = "'0.0929"
$ws.Range("E32").Value = "  +4.07%  "
$ws.Range("D33").Value = "'21.53"
$ws.Range("E33").Value = "  -1.25%  "
$ws.Range("E34").Value = "  +5.59%  "
$ws.Range("E35").Value = "  +1.37%  "
$ws.Range("D36").Value = "'4.67"
$ws.Range("E36").Value = "  +6.82%  "
$ws.Range("D37").Value = "'0.0367"
$ws.Range("E37").Value = "  +4.39%  "
$ws.Range("E38").Value = "  +0.56%  "
$ws.Range("D39").Value = "'3.85"
$ws.Range("E39").Value = "  +15.22%  "
$ws.Range("D40").Value = "'0.245"
$ws.Range("E40").Value = "  +2.22%  "
$ws.Range("E41").Value = "  +3.26%  "
$ws.Range("D42").Value = "'13.64"
$ws.Range("E42").Value = "  +10.57%  "
$ws.Range("D43").Value = "'72.66"
$ws.Range("E43").Value = "  +10.97%  "
$ws.Range("D44").Value = "'6.27"
$ws.Range("E44").Value = "  +21.85%  "
$ws.Range("E45").Value = "  +0.15%  "
$ws.Range("E46").Value = "  +2.36%  "
$ws.Range("E47").Value = "  +0.26%  "
$ws.Range("D48").Value = "'102.26"
$ws.Range("E48").Value = "  +5.75%  "
$ws.Range("D49").Value = "'0.0990"
$ws.Range("E49").Value = "  -2.85%  "
$ws.Range("E50").Value = "  +2.87%  "
$ws.Range("E51").Value = "  +4.64%  "
